$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 5
$ws.Range("A4").Value = 9
$ws.Range("D4").Value = 11
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 3
$ws.Range("D5").Value = 5
$ws.Range("A6").Value = 2
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 10
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 10
$ws.Range("C10").Value = 15
$ws.Range("D10").Value = 15
